$d = $word.ActiveDocument

# 1. Cover-letter date: "2 July 2024" -> "5 July 2024"
$d.Content.Find.Execute("2 July 2024", $true, $false, $false, $false, $false, $true, 1, $false, "5 July 2024", 2) | Out-Null

# 2. Append paper code to first title paragraph
$d.Content.Find.Execute("Part I: Foundations of a rigorous analytical framework", $true, $false, $false, $false, $false, $true, 1, $false, "Part I: Foundations of a rigorous analytical framework (23-007)", 2) | Out-Null

# 3. Append paper code to second title paragraph
$d.Content.Find.Execute("Part II: Applications of the framework", $true, $false, $false, $false, $false, $true, 1, $false, "Part II: Applications of the framework (23-008)", 2) | Out-Null

# 4. "We re-submitted" -> "Following reviewer guidance, we re-submitted"
$p25 = $d.Paragraphs(25).Range
$p25.Find.Execute("We re-submitted", $true, $false, $false, $false, $false, $true, 1, $false, "Following reviewer guidance, we re-submitted", 2) | Out-Null

# 5. Add "(Parts I and II)" before "are being submitted"
$d.Content.Find.Execute("revised versions of both papers are being submitted", $true, $false, $false, $false, $false, $true, 1, $false, "revised versions of both papers (Parts I and II) are being submitted", 2) | Out-Null

# 6. Insert a blank paragraph after the "include:" paragraph (before "All monetary flows...")
$foundRange = $d.Content.Duplicate
$foundRange.Find.Execute("Major improvements to the papers include:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$foundRange.Collapse(0)
$foundRange.InsertParagraphAfter()

# 7. "All monetary flows are now discounted" -> add period
$d.Content.Find.Execute("All monetary flows are now discounted", $true, $false, $false, $false, $false, $true, 1, $false, "All monetary flows are now discounted.", 2) | Out-Null

# 8. "Three utility models are now compared in a new table" -> add colon
$d.Content.Find.Execute("Three utility models are now compared in a new table", $true, $false, $false, $false, $false, $true, 1, $false, "Three utility models are now compared in a new table:", 2) | Out-Null

# 9. "...from the original submission" -> add period
$d.Content.Find.Execute("(CES) utility model from the original submission", $true, $false, $false, $false, $false, $true, 1, $false, "(CES) utility model from the original submission.", 2) | Out-Null

# 10. "...satiated energy service consumption " -> "...satiated energy service consumption."
$d.Content.Find.Execute("energy service consumption ", $true, $false, $false, $false, $false, $true, 1, $false, "energy service consumption.", 2) | Out-Null

# 11. "...a new constant price elasticity (CPE) utility model" -> add period
$d.Content.Find.Execute("constant price elasticity (CPE) utility model", $true, $false, $false, $false, $false, $true, 1, $false, "constant price elasticity (CPE) utility model.", 2) | Out-Null

# 12. "...motivate a value for k" -> add period
$d.Content.Find.Execute("motivate a value for k", $true, $false, $false, $false, $false, $true, 1, $false, "motivate a value for k.", 2) | Out-Null

# 13. "A producer-sided energy price rebound is now estimated" -> add period
$d.Content.Find.Execute("A producer-sided energy price rebound is now estimated", $true, $false, $false, $false, $false, $true, 1, $false, "A producer-sided energy price rebound is now estimated.", 2) | Out-Null

# 14. "this paper" -> "these papers"
$d.Content.Find.Execute("We look forward to seeing this paper in", $true, $false, $false, $false, $false, $true, 1, $false, "We look forward to seeing these papers in", 2) | Out-Null

# 15. Remove the trailing empty paragraph after "Matthew K. Heun"
$d.Paragraphs($d.Paragraphs.Count).Range.Delete() | Out-Null

Write-Output "done"
